$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "N3" = 5
    "O3" = 1.22
    "P3" = 2.38
    "Q3" = 1.65
    "R3" = 1.55
    "S3" = 2.6
    "T3" = 1.58
    "U3" = 2.5
    "X3" = 22
    "Y3" = 16.5
    "Z3" = 24
    "AA3" = 1000
    "AC3" = 9.6
    "AD3" = 13.5
    "AE3" = 34
    "AF3" = 19
    "AG3" = 12
    "AI3" = 36
    "AN3" = 14.5
    "F6" = 1.77
    "G6" = 1.87
    "I6" = 5.4
    "F7" = 5.9
    "G7" = 6
    "T7" = 1.74
    "Q8" = 1.98
    "F12" = 5.2
    "J12" = 3.55
    "K12" = 3.6
    "R12" = 1.25
    "Q13" = 2.28
    "T13" = 2.04
    "F14" = 2.58
    "K14" = 3.3
    "N14" = 2.98
    "P14" = 1.67
    "Q14" = 2.3
    "G16" = 1.51
    "K16" = 4.8
    "N16" = 3.85
    "P16" = 1.99
    "Q16" = 1.89
    "I17" = 4.7
    "J17" = 3.05
    "N17" = 2.74
    "O17" = 1.48
    "P17" = 1.59
    "Q17" = 2.46
    "G18" = 2.06
    "H18" = 4.4
    "I18" = 5.4
    "N18" = 2.92
    "O18" = 1.43
    "Q18" = 2.28
    "G19" = 2.84
    "I19" = 3.15
    "J22" = 3.3
    "F23" = 8.4
    "G23" = 9.6
    "H23" = 1.39
    "I23" = 1.44
    "J23" = 5.1
    "K23" = 5.8
    "N23" = 5.1
    "Z23" = 9.800000000000001
    "AH23" = 25
    "AI23" = 1000
    "AJ23" = 330
    "F24" = 1.82
    "G24" = 1.91
    "H24" = 4.7
    "I24" = 5.3
    "P24" = 1.91
    "I25" = 9
    "J25" = 5.3
    "K25" = 5.4
    "P25" = 2.32
    "H26" = 4.7
    "N26" = 3.85
    "P26" = 1.98
    "Q27" = 1.86
    "S27" = 3.15
    "T27" = 1.69
    "AA27" = 42
    "AE27" = 28
    "P28" = 2.78
    "AK28" = 13.5
    "AN28" = 3.95
    "P30" = 1.72
    "T30" = 2.04
    "I32" = 5.9
    "J32" = 3.35
    "N33" = 3.6
    "Q33" = 1.95
    "O34" = 1.47
    "P34" = 1.66
    "AO34" = 42
    "Q35" = 1.75
    "R35" = 1.5
    "T35" = 1.7
    "U35" = 2.36
    "AA35" = 95
    "AH35" = 17
    "AJ35" = 21
    "AN35" = 11
    "I36" = 8.4
    "Q36" = 1.51
    "P37" = 2.04
    "H39" = 1.5
    "K39" = 4.7
    "AB39" = 25
    "AC39" = 10
    "AJ39" = 410
    "AO39" = 8.6
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Applied $($updates.Count) cell updates"
